$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add hyperlinks in the same order as the original edit (E12, then E3, then E2)
# so the relationship ids line up (rId5 -> E12, rId6 -> E3, rId7 -> E2).
$ws.Hyperlinks.Add($ws.Range("E12"), "mailto:bryanjangeesingh@gmail.com ")
$ws.Range("E12").Value = "bryanjangeesingh@gmail.com "

$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:gtatrinidad@hotmail.com")
$ws.Range("E3").Value = "gtatrinidad@hotmail.com"

$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:jared.hamid@gmail.com")
$ws.Range("E2").Value = "jared.hamid@gmail.com"

# Re-apply the built-in Hyperlink cell style so the new cells share the
# same style record as the pre-existing hyperlink cells (s="1").
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("E3").Style = "Hyperlink"
$ws.Range("E12").Style = "Hyperlink"

# Update the active selection to E4, matching the saved view state.
$ws.Range("E4").Select() | Out-Null
